$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right after the current row 262, shifting all subsequent
# rows (old 263..329) down by one (new 264..330).
$ws.Rows("263:263").Insert(-4121)

# Populate the newly inserted row 263 with the new data entry. The
# categorical columns mirror the row that used to sit there (now row 264),
# while the date/volume/price columns get the new values.
$ws.Range("A263").Value = 5
$ws.Range("B263").Value = "Macroferia Regional de Talca"
$ws.Range("C263").Value = "Maule"
$ws.Range("D263").Value = 44855
$ws.Range("E263").Value = 7
$ws.Range("F263").Value = 100112009
$ws.Range("G263").Value = "Acelga"
$ws.Range("H263").Value = "Sin especificar"
$ws.Range("I263").Value = "Primera"
$ws.Range("J263").Value = 500
$ws.Range("K263").Value = 2500
$ws.Range("L263").Value = 2500
$ws.Range("M263").Value = 2500
$ws.Range("N263").Value = '$/docena de atados (4 kilos)'
$ws.Range("O263").Value = "Región del Maule"
$ws.Range("P263").Value = 625
$ws.Range("Q263").Value = 4
$ws.Range("R263").Value = "Hortaliza"
